$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.8291344720034
$ws.Range("C2").Value = 11.81307979361696
$ws.Range("D2").Value = 5.239497140747471
$ws.Range("E2").Value = 12.73574010039485
$ws.Range("F2").Value = 25.9905995187379
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("L2").Value = 9.876869094851987
$ws.Range("M2").Value = 14.88139831220174
$ws.Range("O2").Value = 23.25365463924796
$ws.Range("B3").Value = 15.28223480662601
$ws.Range("C3").Value = 11.61793548042841
$ws.Range("D3").Value = 5.208347467136297
$ws.Range("E3").Value = 12.78793604549684
$ws.Range("F3").Value = 26.04128224526768
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("L3").Value = 9.884781208495243
$ws.Range("M3").Value = 14.7622182145003
$ws.Range("O3").Value = 23.34896619047711
$ws.Range("B4").Value = 14.93775673891458
$ws.Range("C4").Value = 11.49608816262069
$ws.Range("D4").Value = 5.189065526473319
$ws.Range("E4").Value = 12.8216558513349
$ws.Range("F4").Value = 26.08123352351214
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("L4").Value = 9.890994933423979
$ws.Range("M4").Value = 14.69048682200412
$ws.Range("O4").Value = 23.41422011713684
$ws.Range("B5").Value = 14.79540284180091
$ws.Range("C5").Value = 11.44596846980529
$ws.Range("D5").Value = 5.181172078059125
$ws.Range("E5").Value = 12.8358184592356
$ws.Range("F5").Value = 26.09972578195552
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("L5").Value = 9.893868387801396
$ws.Range("M5").Value = 14.66164327331013
$ws.Range("O5").Value = 23.44249747632652
$ws.Range("B6").Value = 14.77165212672969
$ws.Range("C6").Value = 11.43761930416894
$ws.Range("D6").Value = 5.179859304918526
$ws.Range("E6").Value = 12.83819564477867
$ws.Range("F6").Value = 26.10292970653524
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("L6").Value = 9.894366148080387
$ws.Range("M6").Value = 14.6568779227443
$ws.Range("O6").Value = 23.44729452249911
$ws.Range("B7").Value = 14.93584461976255
$ws.Range("C7").Value = 11.49541405830688
$ws.Range("D7").Value = 5.188959213685567
$ws.Range("E7").Value = 12.82184514486044
$ws.Range("F7").Value = 26.08147397399742
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("L7").Value = 9.891032303371183
$ws.Range("M7").Value = 14.69009622725587
$ws.Range("O7").Value = 23.4145946589984
$ws.Range("B8").Value = 15.64248268655193
$ws.Range("C8").Value = 11.74624091424219
$ws.Range("D8").Value = 5.228791076715888
$ws.Range("E8").Value = 12.75339114463158
$ws.Range("F8").Value = 26.00623676707081
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("L8").Value = 9.879316095015895
$ws.Range("M8").Value = 14.84001958004101
$ws.Range("O8").Value = 23.28511662115813
$ws.Range("B9").Value = 16.95123778871624
$ws.Range("C9").Value = 12.22016117634857
$ws.Range("D9").Value = 5.305519951005748
$ws.Range("E9").Value = 12.6323564086523
$ws.Range("F9").Value = 25.92912536055121
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("L9").Value = 9.867077047695393
$ws.Range("M9").Value = 15.14440715931299
$ws.Range("O9").Value = 23.08493323893276
$ws.Range("B10").Value = 17.85633472017368
$ws.Range("C10").Value = 12.55503191110303
$ws.Range("D10").Value = 5.360859205103387
$ws.Range("E10").Value = 12.55140038061518
$ws.Range("F10").Value = 25.91580649815187
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("L10").Value = 9.864600721555163
$ws.Range("M10").Value = 15.37289493485605
$ws.Range("O10").Value = 22.97101314357237
$ws.Range("B11").Value = 18.25416794598631
$ws.Range("C11").Value = 12.70401346838595
$ws.Range("D11").Value = 5.38576952759386
$ws.Range("E11").Value = 12.51628437560361
$ws.Range("F11").Value = 25.91921473818762
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("L11").Value = 9.864881343355657
$ws.Range("M11").Value = 15.4775734364145
$ws.Range("O11").Value = 22.92646853680093
$ws.Range("B12").Value = 18.40270121067935
$ws.Range("C12").Value = 12.75991017815997
$ws.Range("D12").Value = 5.395161043943748
$ws.Range("E12").Value = 12.50323164513854
$ws.Range("F12").Value = 25.92186908058271
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("L12").Value = 9.865189154930242
$ws.Range("M12").Value = 15.51729146619902
$ws.Range("O12").Value = 22.91065351360665
$ws.Range("B13").Value = 18.37080782884127
$ws.Range("C13").Value = 12.74789549217636
$ws.Range("D13").Value = 5.393140314042954
$ws.Range("E13").Value = 12.50603190859895
$ws.Range("F13").Value = 25.92123674699792
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("L13").Value = 9.865113911865768
$ws.Range("M13").Value = 15.50873442301016
$ws.Range("O13").Value = 22.9140126276652
$ws.Range("B14").Value = 18.26643093408263
$ws.Range("C14").Value = 12.70862273943074
$ws.Range("D14").Value = 5.38654302490515
$ws.Range("E14").Value = 12.51520561719956
$ws.Range("F14").Value = 25.91940577337889
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("L14").Value = 9.864902632742973
$ws.Range("M14").Value = 15.48083964620824
$ws.Range("O14").Value = 22.92514628475451
$ws.Range("B15").Value = 18.20221803122669
$ws.Range("C15").Value = 12.68449833613824
$ws.Range("D15").Value = 5.382496487677439
$ws.Range("E15").Value = 12.52085664303408
$ws.Range("F15").Value = 25.91846189004623
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("L15").Value = 9.864799440521779
$ws.Range("M15").Value = 15.46376270157288
$ws.Range("O15").Value = 22.93210329582386
$ws.Range("B16").Value = 17.83004519201769
$ws.Range("C16").Value = 12.54522491896453
$ws.Range("D16").Value = 5.359225644894494
$ws.Range("E16").Value = 12.55372963180709
$ws.Range("F16").Value = 25.91577449901318
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("L16").Value = 9.864610630798111
$ws.Range("M16").Value = 15.36606644330315
$ws.Range("O16").Value = 22.97407125364056
$ws.Range("B17").Value = 17.59807987380626
$ws.Range("C17").Value = 12.45889931323005
$ws.Range("D17").Value = 5.344879641864935
$ws.Range("E17").Value = 12.57433363314799
$ws.Range("F17").Value = 25.91655276404819
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("L17").Value = 9.864854656765811
$ws.Range("M17").Value = 15.30630238751044
$ws.Range("O17").Value = 23.00168605410845
$ws.Range("B18").Value = 17.46335734665413
$ws.Range("C18").Value = 12.40893362511776
$ws.Range("D18").Value = 5.336603582992676
$ws.Range("E18").Value = 12.58634565303536
$ws.Range("F18").Value = 25.91789143738086
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("L18").Value = 9.865127510826795
$ws.Range("M18").Value = 15.27199882546466
$ws.Range("O18").Value = 23.01825405038185
$ws.Range("B19").Value = 17.41752281989685
$ws.Range("C19").Value = 12.39196342849787
$ws.Range("D19").Value = 5.333797332199373
$ws.Range("E19").Value = 12.59044042810787
$ws.Range("F19").Value = 25.91849762919366
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("L19").Value = 9.865242673938992
$ws.Range("M19").Value = 15.26039728935978
$ws.Range("O19").Value = 23.0239810842265
$ws.Range("B20").Value = 17.62290865891223
$ws.Range("C20").Value = 12.46812154183095
$ws.Range("D20").Value = 5.346409371974207
$ws.Range("E20").Value = 12.57212363262138
$ws.Range("F20").Value = 25.916377675668
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("L20").Value = 9.864814972381357
$ws.Range("M20").Value = 15.31265721047127
$ws.Range("O20").Value = 22.99867550046278
$ws.Range("B21").Value = 18.29714730171771
$ws.Range("C21").Value = 12.72017248123846
$ws.Range("D21").Value = 5.388481962625242
$ws.Range("E21").Value = 12.51250443811944
$ws.Range("F21").Value = 25.91990655238287
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("L21").Value = 9.864959227358318
$ws.Range("M21").Value = 15.48903110611902
$ws.Range("O21").Value = 22.92184743039038
$ws.Range("B22").Value = 18.72541353469969
$ws.Range("C22").Value = 12.8818592054343
$ws.Range("D22").Value = 5.415735415949004
$ws.Range("E22").Value = 12.47496701296451
$ws.Range("F22").Value = 25.93016169913721
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("L22").Value = 9.866227848147263
$ws.Range("M22").Value = 15.6047467754379
$ws.Range("O22").Value = 22.87777660556693
$ws.Range("B23").Value = 18.49800877429065
$ws.Range("C23").Value = 12.7958541815183
$ws.Range("D23").Value = 5.401213173570605
$ws.Range("E23").Value = 12.49487123310018
$ws.Range("F23").Value = 25.92396062076341
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("L23").Value = 9.865443591728344
$ws.Range("M23").Value = 15.54295551341092
$ws.Range("O23").Value = 22.9007341058205
$ws.Range("B24").Value = 17.61168780012354
$ws.Range("C24").Value = 12.46395321586101
$ws.Range("D24").Value = 5.345717869014928
$ws.Range("E24").Value = 12.57312225511605
$ws.Range("F24").Value = 25.91645405713724
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("L24").Value = 9.864832500720434
$ws.Range("M24").Value = 15.30978401974059
$ws.Range("O24").Value = 23.00003441705043
$ws.Range("B25").Value = 16.60647758000777
$ws.Range("C25").Value = 12.09412687057527
$ws.Range("D25").Value = 5.284930391728269
$ws.Range("E25").Value = 12.66369440177342
$ws.Range("F25").Value = 25.94239966343737
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("L25").Value = 9.869241386743608
$ws.Range("M25").Value = 15.06110139386856
$ws.Range("O25").Value = 23.13329408955568
